# Added Week 15 simulations
# Update cumulative season stats on the "Rushing" and "Receiving" sheets.

$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet ---
# Row 2: L.Jackson
$rushing.Range("E2").Value = 30

# Row 3: T.Huntley
$rushing.Range("C3").Value = 6
$rushing.Range("D3").Value = 4
$rushing.Range("E3").Value = 4

# Row 5: L.Murray
$rushing.Range("C5").Value = 43
$rushing.Range("F5").Value = 14

# Row 6: D.Freeman
$rushing.Range("C6").Value = 59
$rushing.Range("D6").Value = 34
$rushing.Range("E6").Value = 9

# Row 10: D.Duvernay
$rushing.Range("D10").Value = 3

# --- Receiving sheet ---
# Row 4: D.Freeman
$receiving.Range("C4").Value = 35
$receiving.Range("D4").Value = 28

# Row 6: M.Brown
$receiving.Range("C6").Value = 69
$receiving.Range("D6").Value = 55
$receiving.Range("E6").Value = 38

# Row 7: S.Watkins
$receiving.Range("C7").Value = 35
$receiving.Range("E7").Value = 13

# Row 8: D.Duvernay
$receiving.Range("C8").Value = 38
$receiving.Range("D8").Value = 27

# Row 11: R.Bateman
$receiving.Range("C11").Value = 34
$receiving.Range("D11").Value = 24
$receiving.Range("E11").Value = 14
$receiving.Range("F11").Value = 8
$receiving.Range("G11").Value = 3
$receiving.Range("H11").Value = 2

# Row 13: M.Andrews
$receiving.Range("C13").Value = 87
$receiving.Range("D13").Value = 61
$receiving.Range("E13").Value = 22
$receiving.Range("F13").Value = 14
$receiving.Range("G13").Value = 16
$receiving.Range("H13").Value = 9

# Row 16: J.Oliver
$receiving.Range("C16").Value = 12
$receiving.Range("D16").Value = 8
